# feat: add 2022-Q4 data
#
# 1) Insert a new "2022-Q4" sheet right after "总计" (pushing 2022-Q3 /
#    2021-Q4 / 2021-Q3 one slot to the right) and populate it with the
#    per-fund holdings table.
# 2) Insert a new summary row at the top of the "总计" sheet's data
#    (row 2) for 2022-Q4, shifting the previous rows down by one.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) New "2022-Q4" detail sheet, placed right after "总计".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $cell = $q4.Cells.Item(1, $i + 2)
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$i]
    $cell.Style = $summary.Cells.Item(1, 2).Style
}

$rows2022Q4 = @(
  ,@(0, "860001", "光大阳光混合A", "9.17", "87.55", "4.01", "0.3677", 8)
  ,@(1, "013396", "华夏新能源车龙头混合C", "8.76", "92.66", "4.08", "0.3574", 8)
  ,@(2, "013395", "华夏新能源车龙头混合A", "6.50", "92.66", "4.08", "0.2652", 8)
  ,@(3, "012846", "恒越蓝筹精选混合", "5.90", "86.65", "2.50", "0.1475", 8)
  ,@(4, "006049", "恒越研究精选混合A/B", "4.23", "88.62", "3.06", "0.1294", 7)
  ,@(5, "011506", "建信高端装备股票A", "3.84", "89.26", "3.14", "0.1206", 10)
  ,@(6, "007192", "恒越研究精选混合C", "3.19", "88.62", "3.06", "0.0976", 7)
  ,@(7, "012924", "华夏新时代灵活配置混合（QDII）美元现汇", "2.09", "77.57", "3.49", "0.0729", 5)
  ,@(8, "012925", "华夏新时代灵活配置混合（QDII）美元现钞", "2.09", "77.57", "3.49", "0.0729", 5)
  ,@(9, "860052", "光大阳光启明星创新驱动主题混合B", "1.73", "87.94", "3.38", "0.0585", 8)
  ,@(10, "860053", "光大阳光启明星创新驱动主题混合C", "1.26", "87.94", "3.38", "0.0426", 8)
  ,@(11, "011507", "建信高端装备股票C", "0.90", "89.26", "3.14", "0.0283", 10)
  ,@(12, "860016", "光大阳光启明星创新驱动主题混合A", "0.45", "87.94", "3.38", "0.0152", 8)
  ,@(13, "860036", "光大阳光混合B", "0.10", "87.55", "4.01", "0.0040", 8)
  ,@(14, "860037", "光大阳光混合C", "0.00", "87.55", "4.01", 0.0, 8)
)

$rowNum = 2
foreach ($r in $rows2022Q4) {
    $q4.Cells.Item($rowNum, 1).Value = $r[0]
    $q4.Cells.Item($rowNum, 1).Style = $summary.Cells.Item(2, 1).Style

    for ($col = 2; $col -le 7; $col++) {
        $val = $r[$col - 1]
        $cell = $q4.Cells.Item($rowNum, $col)
        if ($val -is [string]) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $val
    }

    $q4.Cells.Item($rowNum, 8).Value = $r[7]
    $rowNum++
}

$q4.Range("A1").Select()

# ---------------------------------------------------------------------
# 2) Prepend the 2022-Q4 summary row on "总计".
# ---------------------------------------------------------------------
$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122) # xlPasteFormats

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 15
$summary.Cells.Item(2, 4).Value = 1.78

# Renumber the 0-based index column now that a row was prepended.
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3

$summary.Range("A1").Select()
